$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Three tables (slides 14-16) get a new (built-in) table style id applied.
# ---------------------------------------------------------------------------
for ($i = 14; $i -le 16; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle("{BD9178D3-3440-4623-A82C-790BB8C94A3B}")
        }
    }
}

# ---------------------------------------------------------------------------
# 2) The deck's theme colour scheme (currently "Integral" / Red Violet) is
#    swapped back to the stock "Office" palette - recolour every theme slot
#    via the DrawingML theme colour scheme exposed on a slide.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme

# index -> (name, new RGB as 0xBBGGRR because PowerPoint RGB() is little-endian OLE_COLOR)
$newColors = @(
    0x000000,  # 1  dk1      -> 000000
    0xFFFFFF,  # 2  lt1      -> FFFFFF
    0x6A5444,  # 3  dk2      -> 44546A
    0xE6E6E7,  # 4  lt2      -> E7E6E6
    0xD59B5B,  # 5  accent1  -> 5B9BD5
    0x317DED,  # 6  accent2  -> ED7D31
    0xA5A5A5,  # 7  accent3  -> A5A5A5
    0x00C0FF,  # 8  accent4  -> FFC000
    0xC47244,  # 9  accent5  -> 4472C4
    0x47AD70,  # 10 accent6  -> 70AD47
    0xC16305,  # 11 hlink    -> 0563C1
    0x724F95   # 12 folHlink -> 954F72
)

for ($i = 1; $i -le $newColors.Count; $i++) {
    $themeColor = $themeColors.Colors($i)
    $themeColor.RGB = $newColors[$i - 1]
}
